$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper cell used to force values to be written as literal text (not auto-converted
# to numbers) without mutating any cell style - values are staged here with a Text
# number format, then copied via PasteSpecial(xlPasteValues) onto the real target cell.
$helper = $ws.Range("Z1")
$helper.NumberFormat = "@"

$updates = @(
    @{Cell='D2'; Value='91.145.44'},
    @{Cell='E2'; Value='  +1.29%  '},
    @{Cell='D3'; Value='3.188.30'},
    @{Cell='E3'; Value='  +2.57%  '},
    @{Cell='E4'; Value='  +0.01%  '},
    @{Cell='D5'; Value='220.28'},
    @{Cell='E5'; Value='  +2.63%  '},
    @{Cell='D6'; Value='626.94'},
    @{Cell='E6'; Value='  +0.94%  '},
    @{Cell='D7'; Value='1.06'},
    @{Cell='E7'; Value='  +22.02%  '},
    @{Cell='D8'; Value='0.376'},
    @{Cell='E8'; Value='  +1.05%  '},
    @{Cell='E9'; Value='  -0.02%  '},
    @{Cell='D10'; Value='3.184.19'},
    @{Cell='E10'; Value='  +2.51%  '},
    @{Cell='D11'; Value='0.753'},
    @{Cell='E11'; Value='  +9.47%  '},
    @{Cell='E12'; Value='  +5.45%  '},
    @{Cell='E13'; Value='  +3.66%  '},
    @{Cell='D14'; Value='35.32'},
    @{Cell='E14'; Value='  +8.05%  '},
    @{Cell='D15'; Value='5.59'},
    @{Cell='E15'; Value='  +4.31%  '},
    @{Cell='D16'; Value='90.919.56'},
    @{Cell='E16'; Value='  +0.98%  '},
    @{Cell='E17'; Value='  +1.93%  '},
    @{Cell='D18'; Value='3.177.73'},
    @{Cell='E18'; Value='  +2.23%  '},
    @{Cell='D19'; Value='3.81'},
    @{Cell='E19'; Value='  +11.34%  '},
    @{Cell='D20'; Value='0.0000221'},
    @{Cell='E20'; Value='  +2.23%  '},
    @{Cell='D21'; Value='14.41'},
    @{Cell='E21'; Value='  +5.98%  '},
    @{Cell='D22'; Value='444.58'},
    @{Cell='E22'; Value='  +2.60%  '},
    @{Cell='D23'; Value='9.04'},
    @{Cell='E23'; Value='  +9.16%  '},
    @{Cell='D24'; Value='5.21'},
    @{Cell='E24'; Value='  +3.20%  '},
    @{Cell='E25'; Value='  +9.13%  '},
    @{Cell='B26'; Value='Aptos'},
    @{Cell='C26'; Value='https://coinranking.com/coin/HGYj5JCv5+aptos-apt'},
    @{Cell='D26'; Value='12.44'},
    @{Cell='E26'; Value='  +2.38%  '},
    @{Cell='B27'; Value='Litecoin'},
    @{Cell='C27'; Value='https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'},
    @{Cell='D27'; Value='87.30'},
    @{Cell='E27'; Value='  +1.05%  '},
    @{Cell='D28'; Value='3.345.75'},
    @{Cell='E28'; Value='  +1.76%  '},
    @{Cell='D29'; Value='1.00'},
    @{Cell='E29'; Value='  +0.04%  '},
    @{Cell='E30'; Value='  +1.29%  '},
    @{Cell='D31'; Value='9.36'},
    @{Cell='E31'; Value='  +14.41%  '},
    @{Cell='D32'; Value='0.998'},
    @{Cell='E32'; Value='  -8.67%  '},
    @{Cell='D33'; Value='531.40'},
    @{Cell='E33'; Value='  +3.42%  '},
    @{Cell='B34'; Value='dogwifhat'},
    @{Cell='C34'; Value='https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'},
    @{Cell='D34'; Value='3.78'},
    @{Cell='E34'; Value='  +2.62%  '},
    @{Cell='B35'; Value='EthereumClassic'},
    @{Cell='C35'; Value='https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'},
    @{Cell='D35'; Value='25.09'},
    @{Cell='E35'; Value='  +6.86%  '},
    @{Cell='E36'; Value='  +10.66%  '},
    @{Cell='D37'; Value='7.09'},
    @{Cell='E37'; Value='  +4.79%  '},
    @{Cell='E38'; Value='  +5.69%  '},
    @{Cell='D39'; Value='1.32'},
    @{Cell='E39'; Value='  +5.12%  '},
    @{Cell='B40'; Value='Stellar'},
    @{Cell='C40'; Value='https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'},
    @{Cell='D40'; Value='0.168'},
    @{Cell='E40'; Value='  +17.83%  '},
    @{Cell='B41'; Value='WhiteBITCoin'},
    @{Cell='C41'; Value='https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'},
    @{Cell='D41'; Value='22.23'},
    @{Cell='E41'; Value='  -0.36%  '},
    @{Cell='B42'; Value='FirstDigitalUSD'},
    @{Cell='C42'; Value='https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'},
    @{Cell='D42'; Value='1.00'},
    @{Cell='E42'; Value='  -0.09%  '},
    @{Cell='B43'; Value='PolygonEcosystemToken'},
    @{Cell='C43'; Value='https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'},
    @{Cell='D43'; Value='0.416'},
    @{Cell='E43'; Value='  +11.06%  '},
    @{Cell='B44'; Value='Hedera'},
    @{Cell='C44'; Value='https://coinranking.com/coin/jad286TjB+hedera-hbar'},
    @{Cell='D44'; Value='0.0841'},
    @{Cell='E44'; Value='  +16.44%  '},
    @{Cell='D45'; Value='1.97'},
    @{Cell='E45'; Value='  +5.65%  '},
    @{Cell='E46'; Value='  +0.01%  '},
    @{Cell='D47'; Value='148.55'},
    @{Cell='E47'; Value='  +1.71%  '},
    @{Cell='D48'; Value='1.37'},
    @{Cell='E48'; Value='  +10.75%  '},
    @{Cell='D49'; Value='44.27'},
    @{Cell='E49'; Value='  +1.58%  '},
    @{Cell='D50'; Value='4.41'},
    @{Cell='E50'; Value='  +9.40%  '},
    @{Cell='B51'; Value='ARBITRUM'},
    @{Cell='C51'; Value='https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'},
    @{Cell='D51'; Value='0.656'},
    @{Cell='E51'; Value='  +9.92%  '}
)

foreach ($item in $updates) {
    $helper.Value = $item.Value
    $helper.Copy()
    $ws.Range($item.Cell).PasteSpecial(-4163)
}

$helper.Clear()
$excel.CutCopyMode = 0

Write-Output ("Applied " + $updates.Count + " cell updates")
